$d = $word.ActiveDocument

# In the "Requisitos" bulleted paragraph, the line
#   "LOB1019 -  Física II  (Requisito fraco)"
# needs to move from being the first item to being the last item, i.e.
# after "LOB1004 -  Cálculo II  (Requisito fraco)".
#
# NOTE: this runtime's Range.Text getter is unreliable for zero-length
# (collapsed) ranges, so we avoid reading .Text from collapsed ranges and
# instead rely on Paragraphs.Item(...).Range.Text (non-collapsed) plus
# InsertAfter (which works correctly on collapsed ranges).

$searchText = "LOB1019 -  Física II  (Requisito fraco)"

# Locate the paragraph that contains the line to move.
$paraIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*$searchText*") {
        $paraIndex = $i
        break
    }
}

if ($paraIndex -gt 0) {
    $para = $d.Paragraphs.Item($paraIndex)
    $pStart = $para.Range.Start
    $pEnd = $para.Range.End

    # Search for the target text within the paragraph's range.
    $rng = $d.Range($pStart, $pEnd)
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

    if ($found) {
        # Extend by one character to include the trailing <w:br/> line break.
        [void]$rng.MoveEnd(1, 1)
        $movedText = $rng.Text

        # Remove the text (and its trailing break) from its original spot.
        $rng.Text = ""

        # Re-fetch the paragraph (its End moved after the deletion) and
        # insert the saved text right before the paragraph's ending mark,
        # i.e. after the current last run/break ("Cálculo II" line).
        $para2 = $d.Paragraphs.Item($paraIndex)
        $paraEnd2 = $para2.Range.End
        $insertPoint = $d.Range($paraEnd2 - 1, $paraEnd2 - 1)
        $insertPoint.InsertAfter($movedText)
    }
}
